$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells being updated to Text format so that
# numeric-looking values (e.g. "1.000", "8.650", "136.00") are stored
# verbatim instead of being re-interpreted/truncated as numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.899.20"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "1.741.24"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "248.49"
$ws.Range("E5").Value = "  +5.66%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "0.5123"
$ws.Range("E7").Value = "  -3.24%  "
$ws.Range("D8").Value = "0.2746"
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("D9").Value = "0.06185"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "1.738.44"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").Value = "0.07232"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "15.13"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "0.6484"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "4.629"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "77.65"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "1.000"
$ws.Range("D18").Value = "25.926.17"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").Value = "0.000006813"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "1.964.97"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "4.279"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").Value = "8.650"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("D24").Value = "5.399"
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").Value = "136.00"
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("D26").Value = "1.501"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").Value = "15.24"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "1.774"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").Value = "105.74"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").Value = "3.914"
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").Value = "0.08229"
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("D32").Value = "3.647"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "0.04687"
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("D34").Value = "2.654"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").Value = "0.9988"
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").Value = "0.6256"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "1.918"
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").Value = "1.001"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Value = "100.14"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("D42").Value = "0.7576"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "0.3846"
$ws.Range("E43").Value = "  -2.37%  "
$ws.Range("D44").Value = "4.992"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").Value = "6.298"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").Value = "55.31"
$ws.Range("D48").Value = "0.05231"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").Value = "30.66"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "0.3414"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.485"
$ws.Range("E51").Value = "  -2.35%  "
